# Time_log.docx edit: fix "Use Cases: sistemazione generale" run split,
# and fill in the next (previously empty) row with the new use-case entry.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# --- Hunk 1: merge the two runs "Use Cases: sistemazione " + "generale"
#     into a single run, preserving the existing run formatting. ---
$found = $d.Content.Find.Execute(
    "Use Cases: sistemazione generale", $true, $false, $false, $false,
    $false, $true, 1, $false, "Use Cases: sistemazione generale", 2)

# --- Hunk 2: fill the empty row right after it with the new entry. ---
# Row 46 (1-based) in the single table is the target empty row:
#   col 2/3 (merged) -> Attivita, col 4 -> Data, col 5 -> Tempo (minuti)

# Attivita cell: three runs "Use Cases: sistemazione finale" / " " / "e creazione documento bozza"
$cellA = $tbl.Cell(46, 2)
$posA = $cellA.Range.Start
$a1 = "Use Cases: sistemazione finale"
$a2 = " "
$a3 = "e creazione documento bozza"

$cellA.Range.InsertAfter($a1)
$cellA2 = $tbl.Cell(46, 2)
$cellA2.Range.InsertAfter($a2)
$cellA3 = $tbl.Cell(46, 2)
$cellA3.Range.InsertAfter($a3)

$rA1 = $d.Range($posA, $posA + $a1.Length)
$rA1.Font.Name = "Calibri"
$rA2 = $d.Range($posA + $a1.Length, $posA + $a1.Length + $a2.Length)
$rA2.Font.Name = "Calibri"
$rA3 = $d.Range($posA + $a1.Length + $a2.Length, $posA + $a1.Length + $a2.Length + $a3.Length)
$rA3.Font.Name = "Calibri"

# Data cell: three runs "2" / "6" / "/10/2022"  -> "26/10/2022"
$cellD = $tbl.Cell(46, 4)
$posD = $cellD.Range.Start
$d1 = "2"
$d2 = "6"
$d3 = "/10/2022"

$cellD.Range.InsertAfter($d1)
$cellD2 = $tbl.Cell(46, 4)
$cellD2.Range.InsertAfter($d2)
$cellD3 = $tbl.Cell(46, 4)
$cellD3.Range.InsertAfter($d3)

$rD1 = $d.Range($posD, $posD + $d1.Length)
$rD1.Font.Name = "Calibri"
$rD2 = $d.Range($posD + $d1.Length, $posD + $d1.Length + $d2.Length)
$rD2.Font.Name = "Calibri"
$rD3 = $d.Range($posD + $d1.Length + $d2.Length, $posD + $d1.Length + $d2.Length + $d3.Length)
$rD3.Font.Name = "Calibri"

# Tempo (minuti) cell: single run "120"
$cellT = $tbl.Cell(46, 5)
$posT = $cellT.Range.Start
$t1 = "120"
$cellT.Range.InsertAfter($t1)
$rT1 = $d.Range($posT, $posT + $t1.Length)
$rT1.Font.Name = "Calibri"

Write-Host "Edit applied. Find result:" $found
